$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Combine rows 2-4 (Monk, Token Creature — Monk, Prowess, 1/1) into a single
# python-tuple-like string in A2
$ws.Range("A2").Value = "('Monk', ['Token Creature $([char]0x2014) Monk', 'Prowess', '1/1'])"

# Remove the now-obsolete rows 3-5
$ws.Range("A3:A5").EntireRow.Delete()
